$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at row 7 (pushes existing rows 7-28 down to 8-29) and
#    populate it with the new Crossref AEJ-DOIs data source line.
$ws.Rows(7).Insert()
$ws.Range("A7").Value = 2
$ws.Range("B7").Value = "Crossref (2023)"
$ws.Range("C7").Value = "./data/crossref/crossref_aejdois.Rds"
$ws.Range("D7").Value = "Yes"

# 2. The "mainOA-mapping.xlsx" row (previously row 23, now row 24 after the
#    insertion above) only had Filename/Provided filled in. Fill in its
#    Order and Data Source columns.
$ws.Range("A24").Value = 14
$ws.Range("B24").Value = "Procedural file to map short names to long variable names"

# Keep the active selection consistent with the final layout.
$ws.Range("A25").Select()
